# POE Proposal 1 - edit script
# Summary of changes (per commit message "updated proposal, removed subscription page"):
#  1. Mark the run holding the title picture as <w:noProof/> (cosmetic re-insert marker).
#  2. Remove the "Subscription Service" bullet entirely; the old "Blog Section" bullet
#     becomes the new first bullet and keeps its own text, the old "Blog Section" entry is
#     dropped (its text is now effectively the content of the ex-Subscription paragraph).
#  3. Simplify the "Contact Page" bullet description.
#  4. Insert one extra empty (heading-styled) paragraph after the horizontal rule that
#     precedes "5. Design & User Experience".
#  5. Tidy up a couple of citation paragraphs (merge adjacent runs with identical formatting).
#  6. Mark the built-in "Default Paragraph Font" style as semi-hidden.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Picture run gets <w:noProof/>
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shp = $d.InlineShapes(1)
    $shp.Range.NoProofing = $true
}

# ---------------------------------------------------------------------------
# 2 & 3) Website-pages bullet list: drop "Subscription Service", fold its
#        place into what used to be "Blog Section", drop the stale duplicate,
#        and tighten the "Contact Page" bullet text.
# ---------------------------------------------------------------------------

# Find the paragraph that starts the "Subscription Service" bullet and the one
# that starts "Blog Section" immediately after it, then delete the whole
# "Subscription Service" paragraph (heading run + description run + the
# paragraph mark), leaving "Blog Section" as the entry that follows.
$found = $d.Content.Find.Execute("Subscription Service", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $d.Content.Find.Parent.Paragraphs(1)
}

# Locate paragraph objects directly via the Paragraphs collection so we can
# manipulate whole paragraphs (heading run + description run + paragraph
# mark) reliably.
$subscriptionPara = $null
$blogPara = $null
$contactPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text
    if ($txt -like "Subscription Service*") {
        $subscriptionPara = $para
    } elseif ($txt -like "Blog Section*" -and $subscriptionPara -ne $null -and $blogPara -eq $null) {
        $blogPara = $para
    } elseif ($txt -like "Contact Page*" -and $blogPara -ne $null -and $contactPara -eq $null) {
        $contactPara = $para
    }
}

if ($subscriptionPara -ne $null) {
    # Delete the entire "Subscription Service" paragraph, including its
    # trailing paragraph mark, so what used to be "Blog Section" shifts up
    # and becomes the entry in its place.
    $subscriptionPara.Range.Delete()
}

# Tighten the "Contact Page" description.
$d.Content.Find.Execute("Customer support, FAQs, and a chatbot.", $false, $false, $false, $false, $false, $true, 1, $false, "Customer support.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert an extra empty, heading-styled paragraph right after the
#    horizontal-rule paragraph that precedes "5. Design & User Experience".
# ---------------------------------------------------------------------------
$hrFound = $d.Content.Find.Execute("5. Design & User Experience", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$designHeadingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*5. Design*User Experience*") {
        $designHeadingPara = $para
        break
    }
}

if ($designHeadingPara -ne $null) {
    $insertionRange = $designHeadingPara.Range.Paragraphs(1).Previous().Range
    $newPara = $insertionRange.Paragraphs.Add($insertionRange)
    $newPara.Range.Font.Bold = $true
    $newPara.Range.Font.Size = 14
    $newPara.Range.Font.Color = 4431180
}

# ---------------------------------------------------------------------------
# 5) Citation tidy-ups (merge adjacent identically-formatted runs).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("JBT (no year) Fresh Produce Solutions | ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$d.Content.Find.Execute(". [Online image] Available at: https://www.istockphoto.com/photo/shopping-bag-full-of-fresh-vegetables-and-fruits-gm1128687123-297902712 (Accessed: 6 April 2025).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# ---------------------------------------------------------------------------
# 6) Mark Default Paragraph Font style as semi-hidden.
# ---------------------------------------------------------------------------
try {
    $dpf = $d.Styles("Default Paragraph Font")
    $dpf.Hidden = $true
} catch {
}

Write-Output "done"
